# 3desfreq.xlsx - "Projetos aula02 Lista Mega"
# - Fill in a new "RMST" attendance column (D) on the FREQ sheet for rows 3-20
# - Make FREQ the active sheet/tab (was PONTOS), with A1 selected
# - PONTOS is no longer the active tab; its selection moves to A2

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("FREQ")
$ws2 = $wb.Worksheets.Item("PONTOS")

# New RMST (3rd project) attendance values for FREQ!D3:D20
$ws1.Range("D3").Value  = "P"
$ws1.Range("D4").Value  = "P"
$ws1.Range("D5").Value  = "P"
$ws1.Range("D6").Value  = "P"
$ws1.Range("D7").Value  = "F"
$ws1.Range("D8").Value  = "F"
$ws1.Range("D9").Value  = "F"
$ws1.Range("D10").Value = "F"
$ws1.Range("D11").Value = "P"
$ws1.Range("D12").Value = "F"
$ws1.Range("D13").Value = "P"
$ws1.Range("D14").Value = "F"
$ws1.Range("D15").Value = "P"
$ws1.Range("D16").Value = "P"
$ws1.Range("D17").Value = "P"
$ws1.Range("D18").Value = "P"
$ws1.Range("D19").Value = "P"
$ws1.Range("D20").Value = "P"

# Move the "active tab" / selection state: PONTOS -> A2 (no longer active tab)
$ws2.Range("A2").Select()

# FREQ becomes the active tab, with A1 selected
$ws1.Activate()
$ws1.Range("A1").Select()
